$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This handback-status report regenerated with new source/target file names
# and new timestamps:
#   23f6bd92-b7d1-4908-94ed-2075f84ea54d -> ccbd7795-b258-482a-b40c-ae955711cc8a
#   2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b -> ffff16413676-89f6-4539-a08f-d0c8e8575bd6
#   handoff hash e9ab5f197ccc10b9c3db6e169080ed05ded0cdae / 9324d7989e3346b39dcea40cb541305b6d28540b
#     -> b16176464e09a6896749547b167f4b1fdf4f19ca (shared by both rows now)
# ---------------------------------------------------------------------------

# ---- Sheet "Overview" ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$ov.Range("B2").Value2 = "e2e\ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$ov.Range("G2").Value2 = "2016-08-27 11:02:09"
$ov.Range("A3").Value2 = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
$ov.Range("B3").Value2 = "e2e\ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
$ov.Range("G3").Value2 = "2016-08-27 11:02:09"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\ccbd7795-b258-482a-b40c-ae955711cc8a.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
    }
}

# ---- Sheet "zh-cn" ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$zh.Range("G2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.zh-cn.xlf"
$zh.Range("H2").Value2 = "2016-08-27 11:02:00"
$zh.Range("I2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$zh.Range("J2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.zh-cn.xlf"
$zh.Range("K2").Value2 = "2016-08-27 11:02:29"
$zh.Range("A3").Value2 = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
$zh.Range("G3").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.zh-cn.xlf"
$zh.Range("H3").Value2 = "2016-08-27 11:02:00"
$zh.Range("I3").Value2 = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
$zh.Range("J3").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.zh-cn.xlf"
$zh.Range("K3").Value2 = "2016-08-27 11:02:29"

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
    }
}

# ---- Sheet "de-de" ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$de.Range("G2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.de-de.xlf"
$de.Range("H2").Value2 = "2016-08-27 11:02:09"
$de.Range("I2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$de.Range("J2").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.de-de.xlf"
$de.Range("K2").Value2 = "2016-08-27 11:02:36"
$de.Range("A3").Value2 = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
$de.Range("G3").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.de-de.xlf"
$de.Range("H3").Value2 = "2016-08-27 11:02:09"
$de.Range("I3").Value2 = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
$de.Range("J3").Value2 = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.de-de.xlf"
$de.Range("K3").Value2 = "2016-08-27 11:02:36"

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"
    }
}
